$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking row: correct per-question marks (B11) from 3 to 5
$ws.Range("B11").Value = 5

# Update total row: total marks (B12) from 60 to 100
$ws.Range("B12").Value = 100

# Update correct/total marks summary (E12) from "57/84" to "100/140"
$ws.Range("E12").Value = "100/140"
